# Coverage.xlsx "gsc-export" update
#
# The "Chart" sheet holds a daily coverage export with columns:
#   A: Date, B: Not indexed, C: Indexed, D: Impressions
#
# The refreshed export no longer contains data for the two oldest dates
# (2025-11-08 and 2025-11-09); the row that used to hold 2025-11-10's data
# becomes the new first data row, but its "Not indexed"/"Indexed" counts are
# not yet available (blank) while the Impressions figure for that date is
# unchanged (18).
#
# Net effect: remove the two oldest data rows (rows 2 and 3) and blank out
# the "Not indexed"/"Indexed" values on the (new) first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the 2025-11-08 and 2025-11-09 rows; everything below shifts up by two.
$ws.Rows("2:3").Delete()

# The row that now sits at row 2 (2025-11-10) has no Not indexed / Indexed
# counts yet in this refreshed export; its Impressions value is already
# correct after the shift.
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
